$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = -0.1103015140893087
$ws.Range("C2").Value = 1.867466406311011
$ws.Range("D2").Value = 16.07388988832066
$ws.Range("E2").Value = 4.009225597084886
$ws.Range("F2").Value = 4.097780106789175
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = 0.1241568715908939
$ws.Range("C3").Value = 1.713852921858244
$ws.Range("D3").Value = 10.94170538271163
$ws.Range("E3").Value = 3.307824871832188
$ws.Range("F3").Value = 3.383280951883025
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = -0.5207749427252132
$ws.Range("C4").Value = 1.008444475967709
$ws.Range("D4").Value = 4.236175762135406
$ws.Range("E4").Value = 2.058197211672245
$ws.Range("F4").Value = 2.040396452217854
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = 0.08713849139876084
$ws.Range("C5").Value = 0.7045681055244151
$ws.Range("D5").Value = 1.740827324384067
$ws.Range("E5").Value = 1.319404155057906
$ws.Range("F5").Value = 1.350724642826101
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = 0.05311708210651798
$ws.Range("C6").Value = 0.7222531657823513
$ws.Range("D6").Value = 1.611489849032536
$ws.Range("E6").Value = 1.269444701053392
$ws.Range("F6").Value = 1.30308821368913
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = 0.01086339460369382
$ws.Range("C7").Value = 0.5852842659497787
$ws.Range("D7").Value = 0.6191602100143899
$ws.Range("E7").Value = 0.7868673395270577
$ws.Range("F7").Value = 0.8096026454595719
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = 0.09509580837041771
$ws.Range("C8").Value = 0.5696764806980906
$ws.Range("D8").Value = 0.5776344400534255
$ws.Range("E8").Value = 0.760022657592144
$ws.Range("F8").Value = 0.7772568294901242
$ws.Range("G8").Value = 17

$ws.Range("B9").Value = 0.2285112237303573
$ws.Range("C9").Value = 0.484995177258808
$ws.Range("D9").Value = 0.3855122973504174
$ws.Range("E9").Value = 0.6208963660309322
$ws.Range("F9").Value = 0.5962504332730643
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = 0.2058289907268505
$ws.Range("C10").Value = 0.4532524315174537
$ws.Range("D10").Value = 0.3699707508218096
$ws.Range("E10").Value = 0.6082522098782788
$ws.Range("F10").Value = 0.5924572112923694
$ws.Range("G10").Value = 15

$ws.Range("B11").Value = 0.2381956871565345
$ws.Range("C11").Value = 0.3847249399981075
$ws.Range("D11").Value = 0.1907396565338637
$ws.Range("E11").Value = 0.4367375144567543
$ws.Range("F11").Value = 0.3798820258321563
$ws.Range("G11").Value = 14

